$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$title = $s.Shapes.Item(1)
$tr = $title.TextFrame.TextRange

# Original title text: "Tasks this week"
# Replace the first 11 characters ("Tasks this ") with "Tasks last "
# so the text reads "Tasks last week", split as two runs:
#   "Tasks last " + "week"
$head = $tr.Characters(1, 11)
$head.Text = "Tasks last "
